# The author's edit removed the "pc" data row (old row 23) entirely - deleting
# the whole row shifts every row below it up by one - and appended a brand new
# row of data ("zy_r") at the bottom of the table (new row 48), keeping the
# table's overall extent at A1:D48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Remove the row that held the "pc" entry (row 23). This shifts rows 24:48
# up to become rows 23:47 and automatically drops the now-unused "pc" shared
# string from the workbook.
$ws.Rows("23:23").Delete()

# Append the new "zy_r" record as the new last row of the table (row 48).
$ws.Range("A48").Value = "zy_r"
$ws.Range("B48").Value = -51.39
$ws.Range("C48").Value = -4.05
$ws.Range("D48").Value = 74.78

# Reflect the author's final view/selection state: scrolled down with D48,
# the newly entered cell, selected.
$ws.Range("D48").Select()
